# Atualização de bases das ligas, do dia: 02-03-2024 às 08:34
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 91 and 92: swap match records (id 6924568 / 6924569) ---
# Row 91 (id -> 6924569, Venados FC vs Dorados)
$ws.Range("B91").Value = 6924569
$ws.Range("F91").Value = "Venados FC"
$ws.Range("G91").Value = "Dorados"
$ws.Range("H91").Value = 4
$ws.Range("I91").Value = 1
$ws.Range("J91").Value = "H"
$ws.Range("K91").Value = 1.615
$ws.Range("L91").Value = 4
$ws.Range("M91").Value = 4.5
$ws.Range("N91").Value = 1.5
$ws.Range("O91").Value = 4.75
$ws.Range("P91").Value = 5.75
$ws.Range("Q91").Value = -1.25
$ws.Range("R91").Value = 1.925
$ws.Range("S91").Value = 1.875
$ws.Range("T91").Value = 3
$ws.Range("U91").Value = 1.75
$ws.Range("V91").Value = 1.95
$ws.Range("W91").Value = 0.5
$ws.Range("X91").Value = -1
$ws.Range("Y91").Value = -1
$ws.Range("Z91").Value = 0.925
$ws.Range("AA91").Value = -1
$ws.Range("AB91").Value = 0.75
$ws.Range("AC91").Value = -1

# Row 92 (id -> 6924568, Atletico Morelia vs Atlante)
$ws.Range("B92").Value = 6924568
$ws.Range("F92").Value = "Atletico Morelia"
$ws.Range("G92").Value = "Atlante"
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 1
$ws.Range("J92").Value = "A"
$ws.Range("K92").Value = 2.4
$ws.Range("L92").Value = 3
$ws.Range("M92").Value = 2.875
$ws.Range("N92").Value = 2.7
$ws.Range("O92").Value = 3.1
$ws.Range("P92").Value = 2.8
$ws.Range("Q92").Value = 0
$ws.Range("R92").Value = 1.85
$ws.Range("S92").Value = 1.95
$ws.Range("T92").Value = 2.25
$ws.Range("U92").Value = 1.975
$ws.Range("V92").Value = 1.725
$ws.Range("W92").Value = -1
$ws.Range("X92").Value = -1
$ws.Range("Y92").Value = 1.8
$ws.Range("Z92").Value = -1
$ws.Range("AA92").Value = 0.95
$ws.Range("AB92").Value = -1
$ws.Range("AC92").Value = 0.7250000000000001

# --- Row 176: add result (0-2 A) and updated closing odds ---
$ws.Range("H176").Value = 0
$ws.Range("I176").Value = 2
$ws.Range("J176").Value = "A"
$ws.Range("N176").Value = 6
$ws.Range("O176").Value = 4
$ws.Range("P176").Value = 1.571
$ws.Range("Q176").Value = 1
$ws.Range("R176").Value = 1.85
$ws.Range("S176").Value = 1.95
$ws.Range("T176").Value = 2.5
$ws.Range("U176").Value = 1.975
$ws.Range("V176").Value = 1.825
$ws.Range("W176").Value = -1
$ws.Range("X176").Value = -1
$ws.Range("Y176").Value = 0.571
$ws.Range("Z176").Value = -1
$ws.Range("AA176").Value = 0.95
$ws.Range("AB176").Value = -1
$ws.Range("AC176").Value = 0.825

# --- Row 177: add result (0-1 A) and updated closing odds ---
$ws.Range("H177").Value = 0
$ws.Range("I177").Value = 1
$ws.Range("J177").Value = "A"
$ws.Range("N177").Value = 1.727
$ws.Range("O177").Value = 3.8
$ws.Range("P177").Value = 4.5
$ws.Range("R177").Value = 1.925
$ws.Range("S177").Value = 1.875
$ws.Range("U177").Value = 1.8
$ws.Range("V177").Value = 2
$ws.Range("W177").Value = -1
$ws.Range("X177").Value = -1
$ws.Range("Y177").Value = 3.5
$ws.Range("Z177").Value = -1
$ws.Range("AA177").Value = 0.875
$ws.Range("AB177").Value = -1
$ws.Range("AC177").Value = 1

# --- Row 178: add result (1-1 D) and updated closing odds ---
$ws.Range("H178").Value = 1
$ws.Range("I178").Value = 1
$ws.Range("J178").Value = "D"
$ws.Range("N178").Value = 1.65
$ws.Range("O178").Value = 4
$ws.Range("P178").Value = 5.25
$ws.Range("R178").Value = 1.95
$ws.Range("S178").Value = 1.85
$ws.Range("T178").Value = 2.5
$ws.Range("U178").Value = 1.825
$ws.Range("V178").Value = 1.975
$ws.Range("W178").Value = -1
$ws.Range("X178").Value = 3
$ws.Range("Y178").Value = -1
$ws.Range("Z178").Value = -1
$ws.Range("AA178").Value = 0.8500000000000001
$ws.Range("AB178").Value = -1
$ws.Range("AC178").Value = 0.9750000000000001

# --- Row 179: updated closing odds only (no result yet) ---
$ws.Range("N179").Value = 1.7
$ws.Range("O179").Value = 3.5
$ws.Range("P179").Value = 5.5
$ws.Range("Q179").Value = -0.75
$ws.Range("R179").Value = 1.825
$ws.Range("S179").Value = 1.975
$ws.Range("T179").Value = 2.5
$ws.Range("U179").Value = 2.025
$ws.Range("V179").Value = 1.775
